$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.302.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.087.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  -0.69%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5232"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4416"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09322"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.600"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.892"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.064.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001158"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06656"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.330"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.04%  "

$ws.Range("E22").Value = "  -0.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.354.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.304"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.78"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.513"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "132.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.134"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.666"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.230"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.665"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.858"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02626"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06839"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.350"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6972"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2207"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6814"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.331"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("E46").Value = "  -0.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.374"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.636"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000348"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.208"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.215"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.46%  "
